$d = $word.ActiveDocument

function Insert-RawXml {
    param($Range, $BodyXml)
    $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $BodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $Range.InsertXML($frag) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) Drop the stray "_GoBack" bookmark that currently sits in the empty
#    paragraph right after the "APPMO-SP" version-control table. It will be
#    re-created further below, right where the document was last edited.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) "asignaciones" -> "Asignaciones de rol" (typed as three separate runs,
#    same run formatting) + re-insert the "_GoBack" bookmark right after it.
# ---------------------------------------------------------------------------
$rngFind = $d.Content
$rngFind.Find.Execute("asignaciones") | Out-Null
$target = $d.Range($rngFind.Start, $rngFind.End)

$pPr = '<w:pPr>' +
    '<w:spacing w:line="360" w:lineRule="auto"/>' +
    '<w:jc w:val="center"/>' +
    '<w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' +
    '<w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:color w:val="000000"/><w:szCs w:val="24"/><w:lang w:eastAsia="es-MX"/></w:rPr>' +
    '</w:pPr>'
$rPr = '<w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:color w:val="000000"/><w:szCs w:val="24"/><w:lang w:eastAsia="es-MX"/></w:rPr>'

$body = '<w:p w:rsidR="00706DC8" w:rsidRPr="00E472FF" w:rsidRDefault="009571D8" w:rsidP="007E49C9">' + $pPr +
    '<w:r>' + $rPr + '<w:t>A</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>signaciones</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> de rol</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

Insert-RawXml $target $body

# ---------------------------------------------------------------------------
# 3) " Project " + "Charter" (separate runs, "Charter" wrapped in proofErr
#    spell-check tags) -> a single run " Project Charter" (no proofErr tags).
# ---------------------------------------------------------------------------
$rngFind = $d.Content
$rngFind.Find.Execute(" Project Charter") | Out-Null
$target = $d.Range($rngFind.Start, $rngFind.End)

$pPr = '<w:pPr>' +
    '<w:pStyle w:val="Prrafodelista"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="29"/></w:numPr>' +
    '<w:spacing w:after="0"/>' +
    '<w:jc w:val="left"/>' +
    '<w:rPr><w:sz w:val="16"/></w:rPr>' +
    '</w:pPr>'

$body = '<w:p w:rsidR="00E5480E" w:rsidRPr="00260675" w:rsidRDefault="00E5480E" w:rsidP="00E05C26">' + $pPr +
    '<w:r><w:rPr><w:sz w:val="16"/></w:rPr><w:t xml:space="preserve"> Project Charter</w:t></w:r>' +
    '</w:p>'

Insert-RawXml $target $body

# ---------------------------------------------------------------------------
# 4) " " + "Scopes" + " " + "Statement" (separate runs, "Scopes"/"Statement"
#    wrapped in proofErr spell-check tags) -> a single run " Scopes Statement".
# ---------------------------------------------------------------------------
$rngFind = $d.Content
$rngFind.Find.Execute(" Scopes Statement") | Out-Null
$target = $d.Range($rngFind.Start, $rngFind.End)

$pPr = '<w:pPr>' +
    '<w:pStyle w:val="Prrafodelista"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="30"/></w:numPr>' +
    '<w:spacing w:after="0"/>' +
    '<w:jc w:val="left"/>' +
    '<w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="16"/><w:szCs w:val="24"/></w:rPr>' +
    '</w:pPr>'
$rPr2 = '<w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="16"/><w:szCs w:val="24"/></w:rPr>'

$body = '<w:p w:rsidR="00E5480E" w:rsidRPr="00CA586D" w:rsidRDefault="00E5480E" w:rsidP="00E05C26">' + $pPr +
    '<w:r>' + $rPr2 + '<w:t xml:space="preserve"> Scopes Statement</w:t></w:r>' +
    '</w:p>'

Insert-RawXml $target $body

Write-Host "Edits applied."
